$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3; this pushes the existing row 3 (preproducciongestion
# environment row) down to row 4, carrying its hyperlink reference along (Excel keeps
# the hyperlink anchored at C3 after the insert, which we fix up further below).
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with the same environment data as row 4 (duplicate of the
# "preproducciongestion" row), but with a new NroSiniestro (payment order number) value.
$ws.Range("B3").Value = $ws.Range("B4").Value()
$ws.Range("D3").Value = $ws.Range("D4").Value()
$ws.Range("E3").Value = $ws.Range("E4").Value()
$ws.Range("F3").Value = "'0420172008483"
$ws.Range("G3").Value = "'" + $ws.Range("G4").Value()

# Rebuild the hyperlinks: the original hyperlink should now live on C4 (it moved down
# with its row), and a new hyperlink (same target URL) is added on the new C3.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C4"), "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do")

# Adding a hyperlink via COM nudges the cell onto a freshly-created style; restore both
# cells to the original "Hipervínculo" cell style so no redundant style entries are left.
$ws.Range("C3").Style = "Hipervínculo"
$ws.Range("C4").Style = "Hipervínculo"

$ws.Range("I9").Select()
